# Applies the "Tweak to the ranking algo" edit: appends new log entries for
# Sept 19-20, 2016 (serials 42632/42633) to the 'Logs' sheet, continuing the
# sheet's existing pattern of leaving blank separator rows between entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Cells.Item(509, 1).Value = "Demo"
$ws.Cells.Item(509, 2).Value = 42632
$ws.Cells.Item(509, 3).Value = "1630"
$ws.Cells.Item(509, 4).Value = "SSB"
$ws.Cells.Item(509, 5).Value = "W141"
$ws.Cells.Item(509, 6).Value = "Please give PT staff one HDMI - VGA adaptor for client's laptop"

$ws.Cells.Item(510, 1).Value = "Operator"
$ws.Cells.Item(510, 2).Value = 42632
$ws.Cells.Item(510, 3).Value = "1700"
$ws.Cells.Item(510, 4).Value = "SSB"
$ws.Cells.Item(510, 5).Value = "W141"
$ws.Cells.Item(510, 6).Value = "Operate event between  17:00-17:45"

$ws.Cells.Item(511, 1).Value = "AV Shutdown"
$ws.Cells.Item(511, 2).Value = 42632
$ws.Cells.Item(511, 3).Value = "2000"
$ws.Cells.Item(511, 4).Value = "SSB"
$ws.Cells.Item(511, 5).Value = "W141"
$ws.Cells.Item(511, 6).Value = "Pick up HDMI -VGA adaptor  and return to STEACIE"

$ws.Cells.Item(512, 1).Value = "Demo"
$ws.Cells.Item(512, 2).Value = 42632
$ws.Cells.Item(512, 3).Value = "1900"
$ws.Cells.Item(512, 4).Value = "ACW"
$ws.Cells.Item(512, 5).Value = "304"

$ws.Cells.Item(517, 1).Value = "Pickup Skype Kit"
$ws.Cells.Item(517, 2).Value = 42633
$ws.Cells.Item(517, 3).Value = "1730"
$ws.Cells.Item(517, 4).Value = "HNE"
$ws.Cells.Item(517, 5).Value = "102"
$ws.Cells.Item(517, 6).Value = "Door access code -         5065#              Return skype kit to OSG 1014L"

$ws.Cells.Item(518, 1).Value = "Demo"
$ws.Cells.Item(518, 2).Value = 42633
$ws.Cells.Item(518, 3).Value = "1630"
$ws.Cells.Item(518, 4).Value = "OSG"
$ws.Cells.Item(518, 5).Value = "1001"

$ws.Cells.Item(519, 1).Value = "Demo"
$ws.Cells.Item(519, 2).Value = 42633
$ws.Cells.Item(519, 3).Value = "1630"
$ws.Cells.Item(519, 4).Value = "OSG"
$ws.Cells.Item(519, 5).Value = "1002"

$ws.Cells.Item(520, 1).Value = "Demo"
$ws.Cells.Item(520, 2).Value = 42633
$ws.Cells.Item(520, 3).Value = "1830"
$ws.Cells.Item(520, 4).Value = "OSG"
$ws.Cells.Item(520, 5).Value = "1014G"
$ws.Cells.Item(520, 6).Value = "Video recording via WinMovie  maker -  web cam abd tripod in OSG 1014L"

$ws.Cells.Item(521, 1).Value = "Demo"
$ws.Cells.Item(521, 2).Value = 42633
$ws.Cells.Item(521, 3).Value = "1830"
$ws.Cells.Item(521, 4).Value = "OSG"
$ws.Cells.Item(521, 5).Value = "1014J"
$ws.Cells.Item(521, 6).Value = "Video recording via WinMovie  maker -  web cam abd tripod in OSG 1014L"

$ws.Cells.Item(522, 1).Value = "Demo"
$ws.Cells.Item(522, 2).Value = 42633
$ws.Cells.Item(522, 3).Value = "1830"
$ws.Cells.Item(522, 4).Value = "OSG"
$ws.Cells.Item(522, 5).Value = "1014K"
$ws.Cells.Item(522, 6).Value = "Video recording via WinMovie  maker -  web cam abd tripod in OSG 1014L"

$ws.Cells.Item(523, 1).Value = "Demo"
$ws.Cells.Item(523, 2).Value = 42633
$ws.Cells.Item(523, 3).Value = "1830"
$ws.Cells.Item(523, 4).Value = "OSG"
$ws.Cells.Item(523, 5).Value = "2001"
$ws.Cells.Item(523, 6).Value = "Video recording via WinMovie  maker -  web cam abd tripod in OSG 1014L"

$ws.Cells.Item(524, 1).Value = "Demo"
$ws.Cells.Item(524, 2).Value = 42633
$ws.Cells.Item(524, 3).Value = "1830"
$ws.Cells.Item(524, 4).Value = "OSG"
$ws.Cells.Item(524, 5).Value = "2002"
$ws.Cells.Item(524, 6).Value = "Video recording via WinMovie  maker -  web cam abd tripod in OSG 1014L"

$ws.Cells.Item(525, 1).Value = "Pickup Skype Kit"
$ws.Cells.Item(525, 2).Value = 42633
$ws.Cells.Item(525, 3).Value = "2200"
$ws.Cells.Item(525, 4).Value = "OSG"
$ws.Cells.Item(525, 5).Value = "1014G"
$ws.Cells.Item(525, 6).Value = "Return web cam and tripod to OSG 1014L"

$ws.Cells.Item(526, 1).Value = "Pickup Skype Kit"
$ws.Cells.Item(526, 2).Value = 42633
$ws.Cells.Item(526, 3).Value = "2200"
$ws.Cells.Item(526, 4).Value = "OSG"
$ws.Cells.Item(526, 5).Value = "1014K"
$ws.Cells.Item(526, 6).Value = "Return web cam and tripod to OSG 1014L"

$ws.Cells.Item(527, 1).Value = "Pickup Skype Kit"
$ws.Cells.Item(527, 2).Value = 42633
$ws.Cells.Item(527, 3).Value = "2200"
$ws.Cells.Item(527, 4).Value = "OSG"
$ws.Cells.Item(527, 5).Value = "1014J"
$ws.Cells.Item(527, 6).Value = "Return web cam and tripod to OSG 1014L"

$ws.Cells.Item(528, 1).Value = "Pickup Skype Kit"
$ws.Cells.Item(528, 2).Value = 42633
$ws.Cells.Item(528, 3).Value = "2200"
$ws.Cells.Item(528, 4).Value = "OSG"
$ws.Cells.Item(528, 5).Value = "2001"
$ws.Cells.Item(528, 6).Value = "Return web cam and tripod to OSG 1014L"

$ws.Cells.Item(529, 1).Value = "Pickup Skype Kit"
$ws.Cells.Item(529, 2).Value = 42633
$ws.Cells.Item(529, 3).Value = "2200"
$ws.Cells.Item(529, 4).Value = "OSG"
$ws.Cells.Item(529, 5).Value = "2002"
$ws.Cells.Item(529, 6).Value = "Return web cam and tripod to OSG 1014L"

# Rows whose column-F text wraps to two lines at the current column width
# grow to a 30pt row height (matches the wrapped-text autofit on the other
# entries sharing this same comment further up the log).
$ws.Rows.Item(520).RowHeight = 30
$ws.Rows.Item(521).RowHeight = 30
$ws.Rows.Item(522).RowHeight = 30
$ws.Rows.Item(523).RowHeight = 30
$ws.Rows.Item(524).RowHeight = 30

# Keep the sheet view/selection anchored near the newly appended rows,
# mirroring where Excel leaves the cursor after entering this data.
$ws.Application.Goto($ws.Range("A512"), $true)
$ws.Range("F531").Select()
